# Update market-price-derived Leve profit columns (H-N) across all class tables
# Values refreshed by the scheduled market-data runner; commit: "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 644.75
$ws.Range("I2").Value = 137.4
$ws.Range("K2").Value = 137.4
$ws.Range("M2").Value = -24.40000000000001
$ws.Range("H80").Value = 522.5455
$ws.Range("I80").Value = 554.8
$ws.Range("K80").Value = 1664.4
$ws.Range("M80").Value = -666.3999999999999
$ws.Range("H83").Value = 522.5455
$ws.Range("I83").Value = 554.8
$ws.Range("K83").Value = 4993.2
$ws.Range("M83").Value = -1.199999999999818
$ws.Range("H87").Value = 18200
$ws.Range("J87").Value = 20750
$ws.Range("L87").Value = 20750
$ws.Range("N87").Value = -23246
$ws.Range("H90").Value = 18200
$ws.Range("J90").Value = 20750
$ws.Range("L90").Value = 62250
$ws.Range("N90").Value = -74730
$ws.Range("H132").Value = 1748
$ws.Range("I132").Value = 1426.4286
$ws.Range("K132").Value = 4279.2858
$ws.Range("M132").Value = -1749.2858
$ws.Range("H138").Value = 10367.263
$ws.Range("I138").Value = 9438
$ws.Range("J138").Value = 10429.213
$ws.Range("K138").Value = 28314
$ws.Range("L138").Value = 31287.639
$ws.Range("M138").Value = -23174
$ws.Range("N138").Value = -41567.639

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2598.818
$ws.Range("I2").Value = 2598.818
$ws.Range("K2").Value = 2598.818
$ws.Range("M2").Value = -2485.818
$ws.Range("H32").Value = 24325.691
$ws.Range("I32").Value = 15423.404
$ws.Range("J32").Value = 59934.848
$ws.Range("K32").Value = 15423.404
$ws.Range("L32").Value = 59934.848
$ws.Range("M32").Value = -15136.404
$ws.Range("N32").Value = -60508.848
$ws.Range("H55").Value = 37330
$ws.Range("I55").Value = 37330
$ws.Range("K55").Value = 37330
$ws.Range("M55").Value = -37015
$ws.Range("H116").Value = 2598.818
$ws.Range("I116").Value = 2598.818
$ws.Range("K116").Value = 2598.818
$ws.Range("M116").Value = -304.8180000000002
$ws.Range("H122").Value = 2022.6666
$ws.Range("I122").Value = 1958.375
$ws.Range("K122").Value = 5875.125
$ws.Range("M122").Value = -3425.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2598.818
$ws.Range("I3").Value = 2598.818
$ws.Range("K3").Value = 2598.818
$ws.Range("M3").Value = -2484.818
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 596.3333
$ws.Range("I22").Value = 596.3333
$ws.Range("K22").Value = 596.3333
$ws.Range("M22").Value = -246.3333
$ws.Range("H31").Value = 3132
$ws.Range("I31").Value = 3132
$ws.Range("K31").Value = 3132
$ws.Range("M31").Value = -2837
$ws.Range("H34").Value = 3132
$ws.Range("I34").Value = 3132
$ws.Range("K34").Value = 3132
$ws.Range("M34").Value = -2930
$ws.Range("H99").Value = 4612.273
$ws.Range("I99").Value = 5365.2856
$ws.Range("J99").Value = 3294.5
$ws.Range("K99").Value = 5365.2856
$ws.Range("L99").Value = 3294.5
$ws.Range("M99").Value = -3867.2856
$ws.Range("N99").Value = -6290.5
$ws.Range("H122").Value = 3151.5806
$ws.Range("I122").Value = 2844.9443
$ws.Range("J122").Value = 3576.1538
$ws.Range("K122").Value = 8534.832900000001
$ws.Range("L122").Value = 10728.4614
$ws.Range("M122").Value = -6084.832900000001
$ws.Range("N122").Value = -15628.4614
$ws.Range("H126").Value = 4612.273
$ws.Range("I126").Value = 5365.2856
$ws.Range("J126").Value = 3294.5
$ws.Range("K126").Value = 16095.8568
$ws.Range("L126").Value = 9883.5
$ws.Range("M126").Value = -13625.8568
$ws.Range("N126").Value = -14823.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 240.875
$ws.Range("I2").Value = 366.66666
$ws.Range("J2").Value = 165.4
$ws.Range("K2").Value = 2199.99996
$ws.Range("L2").Value = 992.4000000000001
$ws.Range("M2").Value = -2086.99996
$ws.Range("N2").Value = -1218.4
$ws.Range("H23").Value = 144.63637
$ws.Range("I23").Value = 83
$ws.Range("J23").Value = 218.6
$ws.Range("K23").Value = 249
$ws.Range("L23").Value = 655.8
$ws.Range("M23").Value = -14
$ws.Range("N23").Value = -1125.8
$ws.Range("H55").Value = 40359.6
$ws.Range("I55").Value = 599.3333
$ws.Range("K55").Value = 1797.9999
$ws.Range("M55").Value = -1620.9999
$ws.Range("H86").Value = 1051.8572
$ws.Range("I86").Value = 843.5
$ws.Range("J86").Value = 1067.8846
$ws.Range("K86").Value = 2530.5
$ws.Range("L86").Value = 3203.6538
$ws.Range("M86").Value = -1344.5
$ws.Range("N86").Value = -5575.6538
$ws.Range("H89").Value = 1051.8572
$ws.Range("I89").Value = 843.5
$ws.Range("J89").Value = 1067.8846
$ws.Range("K89").Value = 7591.5
$ws.Range("L89").Value = 9610.9614
$ws.Range("M89").Value = -1663.5
$ws.Range("N89").Value = -21466.9614
$ws.Range("H92").Value = 825
$ws.Range("I92").Value = 825
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2475
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1227
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 1513.8572
$ws.Range("I113").Value = 100
$ws.Range("J113").Value = 1749.5
$ws.Range("K113").Value = 300
$ws.Range("L113").Value = 5248.5
$ws.Range("M113").Value = 1870
$ws.Range("N113").Value = -9588.5
$ws.Range("H121").Value = 901973.2
$ws.Range("J121").Value = 1940.1428
$ws.Range("L121").Value = 5820.428400000001
$ws.Range("N121").Value = -8440.4284
$ws.Range("H122").Value = 2518.923
$ws.Range("I122").Value = 2299.5
$ws.Range("J122").Value = 2530.7837
$ws.Range("K122").Value = 20695.5
$ws.Range("L122").Value = 22777.0533
$ws.Range("M122").Value = -18245.5
$ws.Range("N122").Value = -27677.0533
$ws.Range("H131").Value = 27132.111
$ws.Range("I131").Value = 223237.8
$ws.Range("J131").Value = 7121.3267
$ws.Range("K131").Value = 669713.3999999999
$ws.Range("L131").Value = 21363.9801
$ws.Range("M131").Value = -664673.3999999999
$ws.Range("N131").Value = -31443.9801
$ws.Range("H137").Value = 13605.8
$ws.Range("J137").Value = 14499.75
$ws.Range("L137").Value = 43499.25
$ws.Range("N137").Value = -53699.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 42100.434
$ws.Range("I102").Value = 47691.05
$ws.Range("J102").Value = 4829.6665
$ws.Range("K102").Value = 47691.05
$ws.Range("L102").Value = 4829.6665
$ws.Range("M102").Value = -46069.05
$ws.Range("N102").Value = -8073.6665
$ws.Range("H132").Value = 4856.769
$ws.Range("I132").Value = 5087.6924
$ws.Range("J132").Value = 4394.923
$ws.Range("K132").Value = 15263.0772
$ws.Range("L132").Value = 13184.769
$ws.Range("M132").Value = -12733.0772
$ws.Range("N132").Value = -18244.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9335.3125
$ws.Range("J7").Value = 9540.071
$ws.Range("L7").Value = 9540.071
$ws.Range("N7").Value = -9764.071
$ws.Range("H22").Value = 943.0833
$ws.Range("I22").Value = 938.9091
$ws.Range("J22").Value = 989
$ws.Range("K22").Value = 938.9091
$ws.Range("L22").Value = 989
$ws.Range("M22").Value = -643.9091
$ws.Range("N22").Value = -1579
$ws.Range("H27").Value = 943.0833
$ws.Range("I27").Value = 938.9091
$ws.Range("J27").Value = 989
$ws.Range("K27").Value = 938.9091
$ws.Range("L27").Value = 989
$ws.Range("M27").Value = -831.9091
$ws.Range("N27").Value = -1203
$ws.Range("H46").Value = 2114.2856
$ws.Range("I46").Value = 1450
$ws.Range("J46").Value = 2380
$ws.Range("K46").Value = 1450
$ws.Range("L46").Value = 2380
$ws.Range("M46").Value = -1262
$ws.Range("N46").Value = -2756
$ws.Range("H55").Value = 1099.3334
$ws.Range("I55").Value = 1099.3334
$ws.Range("K55").Value = 1099.3334
$ws.Range("M55").Value = -926.3334
$ws.Range("H100").Value = 3606.8
$ws.Range("I100").Value = 2559.7646
$ws.Range("J100").Value = 5831.75
$ws.Range("K100").Value = 2559.7646
$ws.Range("L100").Value = 5831.75
$ws.Range("M100").Value = -2018.7646
$ws.Range("N100").Value = -6913.75
$ws.Range("H126").Value = 9335.3125
$ws.Range("J126").Value = 9540.071
$ws.Range("L126").Value = 28620.213
$ws.Range("N126").Value = -33560.213
$ws.Range("H132").Value = 4017.5715
$ws.Range("I132").Value = 4180.75
$ws.Range("K132").Value = 12542.25
$ws.Range("M132").Value = -10012.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12000
$ws.Range("I32").Value = 10000
$ws.Range("K32").Value = 10000
$ws.Range("M32").Value = -9683
$ws.Range("H46").Value = 168635.81
$ws.Range("J46").Value = 168635.81
$ws.Range("L46").Value = 168635.81
$ws.Range("N46").Value = -169097.81
$ws.Range("H107").Value = 811.6667
$ws.Range("I107").Value = 811.6667
$ws.Range("K107").Value = 2435.0001
$ws.Range("M107").Value = -515.0001000000002
$ws.Range("H132").Value = 3611.6
$ws.Range("I132").Value = 3485.111
$ws.Range("J132").Value = 4750
$ws.Range("K132").Value = 10455.333
$ws.Range("L132").Value = 14250
$ws.Range("M132").Value = -7925.332999999999
$ws.Range("N132").Value = -19310
$ws.Range("H134").Value = 168635.81
$ws.Range("J134").Value = 168635.81
$ws.Range("L134").Value = 505907.43
$ws.Range("N134").Value = -510977.43
$ws.Range("H136").Value = 3608.1667
$ws.Range("I136").Value = 3726.5
$ws.Range("K136").Value = 11179.5
$ws.Range("M136").Value = -8629.5
